# Update "想去人数" (column F) values across sheets to the newly scraped figures.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1929
$ws1.Range("F3").Value = 812
$ws1.Range("F4").Value = 13551
$ws1.Range("F5").Value = 13370
$ws1.Range("F6").Value = 1036
$ws1.Range("F7").Value = 790
$ws1.Range("F8").Value = 27
$ws1.Range("F11").Value = 13
$ws1.Range("F12").Value = 21
$ws1.Range("F13").Value = 717
$ws1.Range("F14").Value = 2118
$ws1.Range("F15").Value = 46
$ws1.Range("F17").Value = 56
$ws1.Range("F18").Value = 95
$ws1.Range("F20").Value = 416
$ws1.Range("F21").Value = 328
$ws1.Range("F23").Value = 473
$ws1.Range("F24").Value = 795
$ws1.Range("F25").Value = 54

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 880
$ws2.Range("F10").Value = 20

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 80

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1929
$ws4.Range("F5").Value = 812
$ws4.Range("F6").Value = 13551
$ws4.Range("F7").Value = 13370
$ws4.Range("F8").Value = 1036
$ws4.Range("F9").Value = 790
$ws4.Range("F10").Value = 27
$ws4.Range("F13").Value = 13
$ws4.Range("F14").Value = 21
$ws4.Range("F15").Value = 717
$ws4.Range("F18").Value = 2118
$ws4.Range("F19").Value = 46
$ws4.Range("F21").Value = 56
$ws4.Range("F22").Value = 95
$ws4.Range("F26").Value = 80
$ws4.Range("F27").Value = 416
$ws4.Range("F28").Value = 328
$ws4.Range("F30").Value = 473
$ws4.Range("F31").Value = 795
$ws4.Range("F33").Value = 880
$ws4.Range("F35").Value = 20
$ws4.Range("F36").Value = 54
